$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 ("sandesh.achari"): a returned book shifts the remaining issued
# books one column to the left; the now-vacated last slot (H2) becomes
# empty again.
$ws.Range("D2").Value = "Mein Kamf"
$ws.Range("E2").Value = "The Lost Symbol"
$ws.Range("F2").Value = "Open"
$ws.Range("G2").Value = "Deception Point"
$ws.Range("H2").ClearContents()
$ws.Range("H2").Interior.Pattern = -4142

# Row 4 ("vishwajeet.vatharkar"): one book returned (shift left) and two
# new books issued.
$ws.Range("D4").Value = "Open"
$ws.Range("E4").Value = "Inferno"
$ws.Range("F4").Value = "The Lost Symbol"
$ws.Range("G4").ClearContents()
$ws.Range("G4").Interior.Pattern = -4142
